$wb = $excel.ActiveWorkbook

$wsZhanlan = $wb.Worksheets.Item("展览")
# row 3
$wsZhanlan.Cells.Item(3, 6).Value = 8953
# row 4
$wsZhanlan.Cells.Item(4, 6).Value = 1968
# row 5
$wsZhanlan.Cells.Item(5, 6).Value = 6596
# row 6
$wsZhanlan.Cells.Item(6, 6).Value = 173
# row 7
$wsZhanlan.Cells.Item(7, 6).Value = 2127
# row 16
$wsZhanlan.Cells.Item(16, 6).Value = 8836
# row 23
$wsZhanlan.Cells.Item(23, 6).Value = 14
# row 25
$wsZhanlan.Cells.Item(25, 6).Value = 78
# row 28
$wsZhanlan.Cells.Item(28, 6).Value = 1031
# row 30
$wsZhanlan.Cells.Item(30, 6).Value = 63
# row 31
$wsZhanlan.Cells.Item(31, 6).Value = 548
# row 33
$wsZhanlan.Cells.Item(33, 6).Value = 17
# row 35
$wsZhanlan.Cells.Item(35, 6).Value = 2282
# row 37
$wsZhanlan.Cells.Item(37, 6).Value = 531
# row 41
$wsZhanlan.Cells.Item(41, 6).Value = 281
# row 44
$wsZhanlan.Cells.Item(44, 6).Value = 1041
# row 45
$wsZhanlan.Cells.Item(45, 6).Value = 94
# row 47
$wsZhanlan.Cells.Item(47, 6).Value = 10
# row 48
$wsZhanlan.Cells.Item(48, 6).Value = 76
# row 49
$wsZhanlan.Cells.Item(49, 6).Value = 3993

$wsBendi = $wb.Worksheets.Item("本地生活")
# row 2
$wsBendi.Cells.Item(2, 6).Value = 2341
# row 3
$wsBendi.Cells.Item(3, 6).Value = 720
# row 4
$wsBendi.Cells.Item(4, 6).Value = 327

$wsQuanbu = $wb.Worksheets.Item("全部类型")
# row 2
$wsQuanbu.Cells.Item(2, 6).Value = 2341
# row 3
$wsQuanbu.Cells.Item(3, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(3, 2).Value = '2024-09-15'
$wsQuanbu.Cells.Item(3, 3).Value = '杭州·2024CJMF·不止音乐节'
$wsQuanbu.Cells.Item(3, 4).Value = '塘子堰路177号 华数产业园隔壁大草坪'
$wsQuanbu.Cells.Item(3, 5).Value = '2024.09.15 13:00-09.16 21:40'
$wsQuanbu.Cells.Item(3, 6).Value = 406
$wsQuanbu.Cells.Item(3, 7).Value = 168
$wsQuanbu.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90522'
$wsQuanbu.Cells.Item(3, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/3PmG2Bq51723192884141.jpeg'
# row 6
$wsQuanbu.Cells.Item(6, 6).Value = 8953
# row 7
$wsQuanbu.Cells.Item(7, 3).Value = '杭州·木灵动漫 二哈和他的白猫师尊主题餐厅'
$wsQuanbu.Cells.Item(7, 4).Value = '平海路124号 杭州湖滨88'
$wsQuanbu.Cells.Item(7, 5).Value = '2024.09.15 00:00-09.30 23:59'
$wsQuanbu.Cells.Item(7, 6).Value = 327
$wsQuanbu.Cells.Item(7, 7).Value = 10
$wsQuanbu.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91251'
$wsQuanbu.Cells.Item(7, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/wLlo6EFv1724642759732.png'
# row 8
$wsQuanbu.Cells.Item(8, 3).Value = '杭州·浮游猫动漫嘉年华'
$wsQuanbu.Cells.Item(8, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$wsQuanbu.Cells.Item(8, 5).Value = '2024.09.15 09:00-09.16 18:00'
$wsQuanbu.Cells.Item(8, 6).Value = 1968
$wsQuanbu.Cells.Item(8, 7).Value = 68
$wsQuanbu.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88498'
$wsQuanbu.Cells.Item(8, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/qsuFy4iv1719569431608.jpeg'
# row 9
$wsQuanbu.Cells.Item(9, 3).Value = '杭州·理想乡动漫展-同人创作者大会'
$wsQuanbu.Cells.Item(9, 4).Value = '金城路785号B1层 萧山人民奥莱公园'
$wsQuanbu.Cells.Item(9, 5).Value = '2024.09.15 10:00-09.16 17:00'
$wsQuanbu.Cells.Item(9, 6).Value = 6596
$wsQuanbu.Cells.Item(9, 7).Value = 80
$wsQuanbu.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83822'
$wsQuanbu.Cells.Item(9, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/toTii6hH1724653069688.jpeg'
# row 10
$wsQuanbu.Cells.Item(10, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(10, 2).Value = '2024-09-15'
$wsQuanbu.Cells.Item(10, 3).Value = '杭州·西溪银泰 布谷布Goods二次元吃谷嘉年华 免票'
$wsQuanbu.Cells.Item(10, 4).Value = '双龙街588号 西溪银泰城'
$wsQuanbu.Cells.Item(10, 5).Value = '2024.09.15 10:00-09.17 20:00'
$wsQuanbu.Cells.Item(10, 6).Value = 173
$wsQuanbu.Cells.Item(10, 7).Value = 30
$wsQuanbu.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89782'
$wsQuanbu.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/iWnJHkey1721737159663.png'
# row 11
$wsQuanbu.Cells.Item(11, 3).Value = '杭州·Eternal时光国乙only展（日+夜场）'
$wsQuanbu.Cells.Item(11, 4).Value = '创意路1号 中国智谷富春园区'
$wsQuanbu.Cells.Item(11, 5).Value = '2024.09.16 09:30-09.17 17:00'
$wsQuanbu.Cells.Item(11, 6).Value = 2127
$wsQuanbu.Cells.Item(11, 7).Value = 75
$wsQuanbu.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89250'
$wsQuanbu.Cells.Item(11, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/VVNYXGdJ1720966510693.png'
# row 12
$wsQuanbu.Cells.Item(12, 3).Value = '杭州·【中秋月圆·早鸟5折】侯小媛2024邓丽君经典金曲演唱会《月亮代表我的心》'
$wsQuanbu.Cells.Item(12, 4).Value = '湖墅南路136-138号 浙话艺术剧院'
$wsQuanbu.Cells.Item(12, 5).Value = '2024.09.16 19:30-09.16 21:00'
$wsQuanbu.Cells.Item(12, 6).Value = 2
$wsQuanbu.Cells.Item(12, 7).Value = 90
$wsQuanbu.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90160'
$wsQuanbu.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/6BkYOXn31722399525085.jpeg'
# row 13
$wsQuanbu.Cells.Item(13, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(13, 2).Value = '2024-09-16'
$wsQuanbu.Cells.Item(13, 3).Value = '杭州·【中秋限定·早鸟特惠】《笑傲江湖》纪念金庸诞辰百年-经典武侠影视金曲音乐会'
$wsQuanbu.Cells.Item(13, 4).Value = '湖墅南路138号 杭州浙话艺术剧院'
$wsQuanbu.Cells.Item(13, 5).Value = '2024.09.16 14:00-09.16 15:30'
$wsQuanbu.Cells.Item(13, 6).Value = 3
$wsQuanbu.Cells.Item(13, 7).Value = 48
$wsQuanbu.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89413'
$wsQuanbu.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/leCNdQ0S1721199147547.jpeg'
# row 14
$wsQuanbu.Cells.Item(14, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(14, 2).Value = '2024-09-16'
$wsQuanbu.Cells.Item(14, 3).Value = '杭州·第五人格同人only'
$wsQuanbu.Cells.Item(14, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$wsQuanbu.Cells.Item(14, 5).Value = '2024.09.16 10:00-09.16 17:00'
$wsQuanbu.Cells.Item(14, 6).Value = 595
$wsQuanbu.Cells.Item(14, 7).Value = 60
$wsQuanbu.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89550'
$wsQuanbu.Cells.Item(14, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/gFZS33XD1721303396870.jpeg'
# row 15
$wsQuanbu.Cells.Item(15, 3).Value = '临安·兮瑶动漫游戏嘉年华'
$wsQuanbu.Cells.Item(15, 4).Value = '锦北街道临天路2636号 杭州海皇世家酒店'
$wsQuanbu.Cells.Item(15, 5).Value = '2024.09.21 10:00-09.21 17:00'
$wsQuanbu.Cells.Item(15, 6).Value = 22
$wsQuanbu.Cells.Item(15, 7).Value = 45
$wsQuanbu.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91603'
$wsQuanbu.Cells.Item(15, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/GEfxG1fq1726210474075.jpeg'
# row 16
$wsQuanbu.Cells.Item(16, 3).Value = '杭州·《LALALAND爱乐之城》浪漫经典名曲烛光音乐会'
$wsQuanbu.Cells.Item(16, 4).Value = '留泗路东山里22号 大美创意园-2号楼'
$wsQuanbu.Cells.Item(16, 5).Value = '2024.09.21 19:30-09.21 21:00'
$wsQuanbu.Cells.Item(16, 6).Value = 1
$wsQuanbu.Cells.Item(16, 7).Value = 50
$wsQuanbu.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91996'
$wsQuanbu.Cells.Item(16, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/Vfb2As9a1725436079117.png'
# row 17
$wsQuanbu.Cells.Item(17, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(17, 2).Value = '2024-09-21'
$wsQuanbu.Cells.Item(17, 3).Value = '杭州·《天空之城》久石让宫崎骏经典作品音乐会'
$wsQuanbu.Cells.Item(17, 4).Value = '武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）'
$wsQuanbu.Cells.Item(17, 5).Value = '2024.09.21 20:00-09.21 21:30'
$wsQuanbu.Cells.Item(17, 6).Value = 4
$wsQuanbu.Cells.Item(17, 7).Value = 100
$wsQuanbu.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90370'
$wsQuanbu.Cells.Item(17, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/9bBDP3811722912606575.jpeg'
# row 18
$wsQuanbu.Cells.Item(18, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(18, 2).Value = '2024-09-21'
$wsQuanbu.Cells.Item(18, 3).Value = '杭州·银泰百货第一届CYLY动漫游戏嘉年华(免票）'
$wsQuanbu.Cells.Item(18, 4).Value = '启智街515号 滨江银泰'
$wsQuanbu.Cells.Item(18, 5).Value = '2024.09.21 10:00-09.21 17:00'
$wsQuanbu.Cells.Item(18, 6).Value = 69
$wsQuanbu.Cells.Item(18, 7).Value = 20
$wsQuanbu.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91496'
$wsQuanbu.Cells.Item(18, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/7woa53kw1725009838707.png'
# row 19
$wsQuanbu.Cells.Item(19, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(19, 2).Value = '2024-09-28'
$wsQuanbu.Cells.Item(19, 3).Value = '杭州·彩虹社同人ONLY——星鸟之歌'
$wsQuanbu.Cells.Item(19, 4).Value = '丰庆路492号建冠龙禾商务中心A幢 杭州华礼宴国际礼宴中心(龙禾商务中心店)'
$wsQuanbu.Cells.Item(19, 5).Value = '2024.09.28 10:00-09.28 18:00'
$wsQuanbu.Cells.Item(19, 6).Value = 79
$wsQuanbu.Cells.Item(19, 7).Value = 79
$wsQuanbu.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91543'
$wsQuanbu.Cells.Item(19, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/0zc8qiBQ1724912733257.jpeg'
# row 20
$wsQuanbu.Cells.Item(20, 6).Value = 8836
# row 26
$wsQuanbu.Cells.Item(26, 6).Value = 14
# row 27
$wsQuanbu.Cells.Item(27, 3).Value = '杭州·弹丸论破only同人展'
$wsQuanbu.Cells.Item(27, 4).Value = '北干街道萧杭路689号 杭州时尚外滩艺术中心'
$wsQuanbu.Cells.Item(27, 5).Value = '2024.10.02 09:30-10.02 17:00'
$wsQuanbu.Cells.Item(27, 6).Value = 78
$wsQuanbu.Cells.Item(27, 7).Value = 80
$wsQuanbu.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91811'
$wsQuanbu.Cells.Item(27, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/fB9EpBgU1724852399775.png'
# row 28
$wsQuanbu.Cells.Item(28, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(28, 2).Value = '2024-10-03'
$wsQuanbu.Cells.Item(28, 3).Value = '杭州·第二届次元格子动漫展嘉宾内场——赵成晨'
$wsQuanbu.Cells.Item(28, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$wsQuanbu.Cells.Item(28, 5).Value = '2024.10.03 09:30-10.03 17:00'
$wsQuanbu.Cells.Item(28, 6).Value = 200
$wsQuanbu.Cells.Item(28, 7).Value = 238
$wsQuanbu.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91121'
$wsQuanbu.Cells.Item(28, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/ddmmhJEE1724313674505.jpeg'
# row 29
$wsQuanbu.Cells.Item(29, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(29, 2).Value = '2024-10-04'
$wsQuanbu.Cells.Item(29, 3).Value = '杭州·创世次元第五人格同人only展'
$wsQuanbu.Cells.Item(29, 4).Value = '小河路与桥弄街交叉口东北50米 桥西历史文化街区'
$wsQuanbu.Cells.Item(29, 5).Value = '2024.10.04 10:00-10.05 17:00'
$wsQuanbu.Cells.Item(29, 6).Value = 1031
$wsQuanbu.Cells.Item(29, 7).Value = 75
$wsQuanbu.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92141'
$wsQuanbu.Cells.Item(29, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/MMF3dkAw1725550270634.jpeg'
# row 30
$wsQuanbu.Cells.Item(30, 3).Value = '杭州·逐月节·园游会·原神×绝区零×崩铁×崩坏同人only'
$wsQuanbu.Cells.Item(30, 4).Value = '莫干山路987号 资辉壹方汇'
$wsQuanbu.Cells.Item(30, 5).Value = '2024.10.04 09:30-10.05 17:00'
$wsQuanbu.Cells.Item(30, 6).Value = 8
$wsQuanbu.Cells.Item(30, 7).Value = 58
$wsQuanbu.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92406'
$wsQuanbu.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/mQh43oPd1726134932363.png'
# row 31
$wsQuanbu.Cells.Item(31, 3).Value = '杭州·2024·华彩的摔跤宴 - 木吉KAZUYA降临'
$wsQuanbu.Cells.Item(31, 4).Value = '莫干山路188-200号 之江饭店(莫干山路店)'
$wsQuanbu.Cells.Item(31, 5).Value = '2024.10.05 10:00-10.05 16:00'
$wsQuanbu.Cells.Item(31, 6).Value = 63
$wsQuanbu.Cells.Item(31, 7).Value = 88
$wsQuanbu.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92402'
$wsQuanbu.Cells.Item(31, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/ZylQGk1P1726033298213.png'
# row 32
$wsQuanbu.Cells.Item(32, 3).Value = '杭州·德山秀典杭州粉丝见面会'
$wsQuanbu.Cells.Item(32, 4).Value = '教工路198号浙商大创业园A幢3楼 杭州子墨汇演中心'
$wsQuanbu.Cells.Item(32, 5).Value = '2024.10.05 12:30-10.05 17:30'
$wsQuanbu.Cells.Item(32, 6).Value = 2
$wsQuanbu.Cells.Item(32, 7).Value = 380
$wsQuanbu.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92385'
$wsQuanbu.Cells.Item(32, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/vqEQxc081726130357271.jpeg'
# row 33
$wsQuanbu.Cells.Item(33, 3).Value = '杭州·文豪野犬同人only2.0'
$wsQuanbu.Cells.Item(33, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$wsQuanbu.Cells.Item(33, 6).Value = 548
$wsQuanbu.Cells.Item(33, 7).Value = 60
$wsQuanbu.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92226'
$wsQuanbu.Cells.Item(33, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/nkCZvaiO1725872765608.jpeg'
# row 34
$wsQuanbu.Cells.Item(34, 3).Value = '杭州·次元幻想【玩美大舞台&吃/换谷大会】（免费活动）'
$wsQuanbu.Cells.Item(34, 4).Value = '文三路 玩美的一天沉浸式生活街区'
$wsQuanbu.Cells.Item(34, 6).Value = 27
$wsQuanbu.Cells.Item(34, 7).Value = 30
$wsQuanbu.Cells.Item(34, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92028'
$wsQuanbu.Cells.Item(34, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/FaEB96xH1725394323651.jpeg'
# row 35
$wsQuanbu.Cells.Item(35, 3).Value = '杭州·火影同人ONLY'
$wsQuanbu.Cells.Item(35, 4).Value = '金城路785号B1层 萧山人民奥莱公园'
$wsQuanbu.Cells.Item(35, 5).Value = '2024.10.05 10:00-10.05 18:00'
$wsQuanbu.Cells.Item(35, 6).Value = 17
$wsQuanbu.Cells.Item(35, 7).Value = 78
$wsQuanbu.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92458'
$wsQuanbu.Cells.Item(35, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/oKXQv6sL1726196644046.png'
# row 36
$wsQuanbu.Cells.Item(36, 3).Value = '杭州·第五人格同人only2.0'
$wsQuanbu.Cells.Item(36, 5).Value = '2024.10.05 10:00-10.05 17:00'
$wsQuanbu.Cells.Item(36, 6).Value = 537
$wsQuanbu.Cells.Item(36, 7).Value = 60
$wsQuanbu.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92229'
$wsQuanbu.Cells.Item(36, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/KGPiiH2U1725873923225.jpeg'
# row 37
$wsQuanbu.Cells.Item(37, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(37, 2).Value = '2024-10-05'
$wsQuanbu.Cells.Item(37, 3).Value = '杭州·首届CCPC动漫嘉年华'
$wsQuanbu.Cells.Item(37, 4).Value = '长乐路29号五组2幢 杭州运河文化发布中心'
$wsQuanbu.Cells.Item(37, 5).Value = '2024.10.05 09:00-10.06 18:00'
$wsQuanbu.Cells.Item(37, 6).Value = 2282
$wsQuanbu.Cells.Item(37, 7).Value = 49.9
$wsQuanbu.Cells.Item(37, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91102'
$wsQuanbu.Cells.Item(37, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/ErwKRZSH1724749999253.jpeg'
# row 38
$wsQuanbu.Cells.Item(38, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(38, 2).Value = '2024-10-05'
$wsQuanbu.Cells.Item(38, 3).Value = '杭州·鸢飞鱼跃代号鸢only'
$wsQuanbu.Cells.Item(38, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$wsQuanbu.Cells.Item(38, 5).Value = '2024.10.05 09:30-10.05 17:00'
$wsQuanbu.Cells.Item(38, 6).Value = 872
$wsQuanbu.Cells.Item(38, 7).Value = 85
$wsQuanbu.Cells.Item(38, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88452'
$wsQuanbu.Cells.Item(38, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/etOXBCrl1719678030944.jpeg'
# row 39
$wsQuanbu.Cells.Item(39, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(39, 2).Value = '2024-10-06'
$wsQuanbu.Cells.Item(39, 3).Value = '杭州·东方乐典2024'
$wsQuanbu.Cells.Item(39, 4).Value = '万塘路262号万塘汇城市生活广场南楼2层 酒球会'
$wsQuanbu.Cells.Item(39, 5).Value = '2024.10.06 18:00-10.06 21:00'
$wsQuanbu.Cells.Item(39, 6).Value = 9
$wsQuanbu.Cells.Item(39, 7).Value = 199
$wsQuanbu.Cells.Item(39, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92460'
$wsQuanbu.Cells.Item(39, 9).Value = '//i1.hdslb.com/bfs/openplatform/202409/k4oqHAK81725949095139.jpeg'
# row 40
$wsQuanbu.Cells.Item(40, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(40, 2).Value = '2024-10-18'
$wsQuanbu.Cells.Item(40, 3).Value = '杭州·Luca Stricagnoli 2024《进化时间》指弹吉他音乐会'
$wsQuanbu.Cells.Item(40, 4).Value = '萧山区杭州国际博览中心A座6楼 九莱福音乐现场'
$wsQuanbu.Cells.Item(40, 5).Value = '2024.10.18 20:00-10.18 21:30'
$wsQuanbu.Cells.Item(40, 6).Value = 6
$wsQuanbu.Cells.Item(40, 7).Value = 220
$wsQuanbu.Cells.Item(40, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91357'
$wsQuanbu.Cells.Item(40, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/0jkK8abH1724743622694.jpeg'
# row 41
$wsQuanbu.Cells.Item(41, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(41, 2).Value = '2024-10-19'
$wsQuanbu.Cells.Item(41, 3).Value = '杭州·SK怀旧动漫展SK12'
$wsQuanbu.Cells.Item(41, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$wsQuanbu.Cells.Item(41, 5).Value = '2024.10.19 09:00-10.20 18:00'
$wsQuanbu.Cells.Item(41, 6).Value = 531
$wsQuanbu.Cells.Item(41, 7).Value = 65
$wsQuanbu.Cells.Item(41, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90921'
$wsQuanbu.Cells.Item(41, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/PspqQq6H1723894652098.jpeg'
# row 42
$wsQuanbu.Cells.Item(42, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(42, 2).Value = '2024-10-26'
$wsQuanbu.Cells.Item(42, 3).Value = '杭州·亿万心动国乙✘代号鸢同人only(日夜场）'
$wsQuanbu.Cells.Item(42, 4).Value = '皓月路299号 诺丁山艺术中心'
$wsQuanbu.Cells.Item(42, 5).Value = '2024.10.26 10:00-10.26 21:00'
$wsQuanbu.Cells.Item(42, 6).Value = 281
$wsQuanbu.Cells.Item(42, 7).Value = 9.9
$wsQuanbu.Cells.Item(42, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91962'
$wsQuanbu.Cells.Item(42, 9).Value = '//i2.hdslb.com/bfs/openplatform/202409/LU32zDTR1725617506119.jpeg'
# row 43
$wsQuanbu.Cells.Item(43, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(43, 2).Value = '2024-11-02'
$wsQuanbu.Cells.Item(43, 3).Value = '杭州·BanGDream! Only同人展'
$wsQuanbu.Cells.Item(43, 4).Value = '石祥路与丽水北路交叉口 大运河音乐公园'
$wsQuanbu.Cells.Item(43, 5).Value = '2024.11.02 10:00-11.03 20:00'
$wsQuanbu.Cells.Item(43, 6).Value = 177
$wsQuanbu.Cells.Item(43, 7).Value = 89
$wsQuanbu.Cells.Item(43, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91168'
$wsQuanbu.Cells.Item(43, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/0vTxEVyz1724222524879.jpeg'
# row 44
$wsQuanbu.Cells.Item(44, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(44, 2).Value = '2024-11-10'
$wsQuanbu.Cells.Item(44, 3).Value = '杭州·崩坏同人ONLY 爱莉希雅生日会'
$wsQuanbu.Cells.Item(44, 4).Value = '康候圣街99号 顺丰创新中心'
$wsQuanbu.Cells.Item(44, 5).Value = '2024.11.10 08:00-11.10 20:00'
$wsQuanbu.Cells.Item(44, 6).Value = 76
$wsQuanbu.Cells.Item(44, 7).Value = 79
$wsQuanbu.Cells.Item(44, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=92228'
$wsQuanbu.Cells.Item(44, 9).Value = '//i0.hdslb.com/bfs/openplatform/202409/1FsO31h71725897488610.jpeg'
# row 45
$wsQuanbu.Cells.Item(45, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(45, 2).Value = '2024-12-07'
$wsQuanbu.Cells.Item(45, 3).Value = '杭州·原神X崩坏X星铁旅行盛宴·同人only首展'
$wsQuanbu.Cells.Item(45, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$wsQuanbu.Cells.Item(45, 5).Value = '2024.12.07 10:00-12.08 17:00'
$wsQuanbu.Cells.Item(45, 6).Value = 3993
$wsQuanbu.Cells.Item(45, 7).Value = 65
$wsQuanbu.Cells.Item(45, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88429'
$wsQuanbu.Cells.Item(45, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/tmiou8M61722221207521.jpeg'
# row 48
$wsQuanbu.Cells.Item(48, 2).NumberFormat = "@"
$wsQuanbu.Cells.Item(48, 2).Value = '2024-12-31'
$wsQuanbu.Cells.Item(48, 3).Value = '杭州·2025大剧院缤纷跨年夜 爱·大声告白-成都“知更”室内合唱团音乐会'
$wsQuanbu.Cells.Item(48, 4).Value = '新业路39号 杭州大剧院'
$wsQuanbu.Cells.Item(48, 5).Value = '2024.12.31 22:30-2025.01.01 00:00'
$wsQuanbu.Cells.Item(48, 6).Value = 2
$wsQuanbu.Cells.Item(48, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91210'
$wsQuanbu.Cells.Item(48, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/RGm2uKFJ1724395472501.jpeg'

